# Doing Updates for Financials
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PAYX")

$ws.Range("D43").Value = 524600
$ws.Range("D45").Value = 4944800
$ws.Range("D46").Value = 5893600
$ws.Range("D52").Value = 377400
$ws.Range("D54").Value = 7915400
$ws.Range("D57").Value = 73700
$ws.Range("D59").Value = 5223000
$ws.Range("D60").Value = 5296700
$ws.Range("D62").Value = 261900
$ws.Range("D66").Value = 5558600
$ws.Range("D72").Value = 1262600
$ws.Range("D76").Value = 2356800

$wb.Save()
